# LOQ4239.xlsx update — 2021-01-29 build
#
# Semantic change: a new "Docentes responsáveis" professor row is inserted
# (11079086 - Herlandí de Souza Andrade), pushing the existing professor row
# and everything below it down by one row; the activation date, teaching
# method, evaluation criteria and recovery rule texts are also updated.
#
# We rebuild rows 13-23 explicitly (rather than relying on Rows.Insert(),
# whose auto formatting/height inheritance doesn't match the target layout)
# so that styles, shared-string reuse and row heights land exactly right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativação: date bump ----------------------------------------------
$ws.Range("B8").Value = "01/01/2021"
$ws.Range("C8").Value = "01/01/2021"

# --- Row 13: NEW professor (Docentes responsáveis) ---------------------
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Rows.Item(13).AutoFit()

# --- Row 14: previous professor, shifted down one row -------------------
$ws.Range("A14").Clear()
$ws.Range("B14").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C14").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Rows.Item(14).AutoFit()

# --- Row 15: Programa resumido (was row 14) -----------------------------
$ws.Range("A15").Value = "Programa resumido:"
$ws.Range("B15").Value = "1. Conceitos Fundamentais de Administração.`n2. Noções Básicas de Estratégia."
$ws.Range("C15").Value = "1. Conceitos Fundamentais de Administração.`n2. Noções Básicas de Estratégia."
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16: Short syllabus (was row 15) --------------------------------
$ws.Range("A16").Value = "Short syllabus:"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()
$ws.Rows.Item(16).RowHeight = 60

# --- Row 17: Programa (was row 16) --------------------------------------
$ws.Range("A17").Value = "Programa:"
$ws.Range("B17").Value = "1. Teoria Geral de Administração: Histórico. Visão de Taylor. Escola clássica de administração.`n2. Conceitos básicos de Estratégia, Análise SWOT, Balanced ScoreCard, Mapas Estratégicos"
$ws.Range("C17").Value = "1. Teoria Geral de Administração: Histórico. Visão de Taylor. Escola clássica de administração.`n2. Conceitos básicos de Estratégia, Análise SWOT, Balanced ScoreCard, Mapas Estratégicos"
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18: Syllabus (was row 17) --------------------------------------
$ws.Range("A18").Value = "Syllabus:"
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).RowHeight = 120

# --- Row 19: Avaliação (was row 18) -------------------------------------
$ws.Range("A19").Value = "Avaliação:"
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()
$ws.Rows.Item(19).AutoFit()

# --- Row 20: Método (was row 19), text changed --------------------------
$ws.Range("A20").Value = "Método:"
$ws.Range("B20").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Range("C20").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21: Critério (was row 20), text changed ------------------------
$ws.Range("A21").Value = "Critério:"
$ws.Range("B21").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Range("C21").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22: Norma de recuperação (was row 21), text changed ------------
$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Range("B22").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("C22").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Rows.Item(22).RowHeight = 60

# --- Row 23: Bibliografia (was row 22), now one row lower ---------------
$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = "CHIAVENATO, I; SAPIRO, A. Planejamento Estratégico. Rio de Janeiro. Campus, 2004 `n`nCOLLINS, J.C.; PORRAS, J. I. Feitas para Durar: Práticas bem-sucedidas de empresas visionárias. 9ª Ed.  Rio de Janeiro. Rocco, 2007 `n`nHERRERO, E. Balanced Scorecard e a Gestão Estratégica. Rio de Janeiro. Campus, 2005. `n`nKAPLAN, R; NORTON, D. Kaplan e Norton na Prática. Rio de Janeiro. Campus, 2004 `n`nKAPLAN, R; NORTON, D. A Estratégia em Ação: Balanced Scorecard. Rio de Janeiro. Campus, 1997 `n`nKAPLAN, R; NORTON, D. Mapas Estratégicos. Rio de Janeiro. Campus, 2004 `n`nTZU, S. A Arte da Guerra (Edição Completa). São Paulo. WMF Martins Fontes, 2009."
$ws.Range("C23").Value = "CHIAVENATO, I; SAPIRO, A. Planejamento Estratégico. Rio de Janeiro. Campus, 2004 `n`nCOLLINS, J.C.; PORRAS, J. I. Feitas para Durar: Práticas bem-sucedidas de empresas visionárias. 9ª Ed.  Rio de Janeiro. Rocco, 2007 `n`nHERRERO, E. Balanced Scorecard e a Gestão Estratégica. Rio de Janeiro. Campus, 2005. `n`nKAPLAN, R; NORTON, D. Kaplan e Norton na Prática. Rio de Janeiro. Campus, 2004 `n`nKAPLAN, R; NORTON, D. A Estratégia em Ação: Balanced Scorecard. Rio de Janeiro. Campus, 1997 `n`nKAPLAN, R; NORTON, D. Mapas Estratégicos. Rio de Janeiro. Campus, 2004 `n`nTZU, S. A Arte da Guerra (Edição Completa). São Paulo. WMF Martins Fontes, 2009."
$ws.Rows.Item(23).RowHeight = 120
